$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 14

# Copy formatting (style) from the row above for column A (date style)
$ws.Cells.Item($row - 1, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item($row, 1).Value = 42619.89266203704
$ws.Cells.Item($row, 2).Value = 26
$ws.Cells.Item($row, 3).Value = 63
$ws.Cells.Item($row, 4).Value = 34
$ws.Cells.Item($row, 5).Value = 63
$ws.Cells.Item($row, 6).Value = 32
$ws.Cells.Item($row, 7).Value = 13111
$ws.Cells.Item($row, 8).Value = 25133
$ws.Cells.Item($row, 9).Value = 2818
$ws.Cells.Item($row, 10).Value = 426
$ws.Cells.Item($row, 11).Value = 230
$ws.Cells.Item($row, 12).Value = 34
$ws.Cells.Item($row, 13).Value = 16
$ws.Cells.Item($row, 14).Value = "Bag"
